# [Feat 2269][Feat 3238] Fixed : parser only processed mandatory columns
#
# The TEST_CASES sheet parser used to only read the "mandatory" columns.
# This test fixture is updated so that the optional TC_REFERENCE (column G)
# and TC_NAME (column H) columns are also populated for each of the 10
# existing test-case rows, exercising the parser fix.

$wb = $excel.ActiveWorkbook

$testCases = $wb.Worksheets.Item(1)   # TEST_CASES
$steps     = $wb.Worksheets.Item(2)   # STEPS

$refs  = @("ref1","ref2","ref3","ref4","ref5","ref6","ref7","ref8","ref9","ref10")
$names = @("name1","name2","name3","name4","name5","name6","name7","name8","name9","name10")

# Fill in the TC_REFERENCE column (G) for all data rows first ...
for ($i = 0; $i -lt $refs.Length; $i++) {
    $row = 2 + $i
    $testCases.Cells.Item($row, 7).Value = $refs[$i]
}

# ... then the TC_NAME column (H) for all data rows.
for ($i = 0; $i -lt $names.Length; $i++) {
    $row = 2 + $i
    $testCases.Cells.Item($row, 8).Value = $names[$i]
}

# Re-apply the banded border formatting across the newly-filled G:H cells so
# it matches the rest of the row (only formats change here, not values).
$noTopBorderSource  = $testCases.Range("D4")   # thin box, no top border
$fullBoxBorderSource = $testCases.Range("A3")  # full thin box border

$rowsNeedingNoTopBorder  = @(4, 6, 8, 10)
$rowsNeedingFullBoxBorder = @(9)

foreach ($row in $rowsNeedingNoTopBorder) {
    $noTopBorderSource.Copy()
    $testCases.Range("G$row`:H$row").PasteSpecial(-4122)
}

foreach ($row in $rowsNeedingFullBoxBorder) {
    $fullBoxBorderSource.Copy()
    $testCases.Range("G$row`:H$row").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# Restore the on-screen selections to what they were left at after the edit.
$steps.Range("A2").Select() | Out-Null
$testCases.Range("H2:H11").Select() | Out-Null
